# Doing Updates for Financials
# Insert a new "most recent period" column before column D on the BMNM
# sheet and shift the existing D:K data right to E:L, then populate the
# new column D with the latest period's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before D - this shifts D:K -> E:L automatically,
#    carrying over values, number formats and styles.
$ws.Columns("D").Insert()

# 2) The freshly inserted column D cells currently inherit column C's
#    format. Copy the (now shifted) formats from column E back onto column
#    D for each of the three data blocks so the new column matches the
#    rest of the table (date format on the header rows, number format on
#    the value rows).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 3) Populate the new column D with the latest period's values.
$numericValues = @{
    7 = 43465
    8 = 6300
    9 = 5500
    10 = 800
    13 = 0
    15 = 0
    17 = 12000
    18 = -5700
    20 = 0
    21 = -5600
    22 = 0
    23 = -5700
    24 = 21100
    25 = 0
    26 = -26800
    27 = -26800
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = -26800
    34 = 0
    35 = -26800
    38 = 43465
    41 = 4900
    42 = 0
    43 = 800
    46 = 0
    47 = 212400
    48 = 3300
    49 = 0
    50 = 0
    51 = 0
    52 = 24500
    53 = 0
    54 = 259400
    57 = 0
    58 = 200400
    59 = 700
    60 = 0
    61 = 26800
    62 = 0
    63 = 0
    64 = 0
    65 = 0
    66 = 230400
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = -306000
    73 = 0
    74 = 0
    75 = 0
    76 = 29000
    77 = 0
    80 = 43465
    81 = -26800
    83 = 100
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 6400
    91 = 0
    92 = 0
    93 = 0
    94 = -9100
    96 = 0
    97 = 0
    98 = 0
    99 = 0
    100 = 200
    101 = 0
    102 = -2500
}

$naRows = @(12, 14, 44, 45)

foreach ($r in $numericValues.Keys) {
    $ws.Cells.Item($r, 4).Value2 = $numericValues[$r]
}

foreach ($r in $naRows) {
    $ws.Cells.Item($r, 4).Value2 = "NA"
}
